# Apply the "July 15, 2020 07:37:50 PM" covid_disparities_output run results
# to the existing workbook (Sheet1), per the provided OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 (California - Los Angeles) ---
$ws.Range("B16").Value = 44026
$ws.Range("C16").Value = 143009
$ws.Range("D16").Value = 3936
$ws.Range("E16").Value = 3818
$ws.Range("F16").Value = 393
$ws.Range("G16").Value = 4.7
$ws.Range("H16").Value = 10.74
$ws.Range("K16").Value = 81305
$ws.Range("L16").Value = 3658

# --- Row 27 (California) ---
$ws.Range("D27").Value = 7227
$ws.Range("G27").Value = 4.34
$ws.Range("H27").Value = 8.82
$ws.Range("J27").Value = $false

# --- Row 36 (Iowa) ---
$ws.Range("C36").Value = 36324
